$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: " her window" / bookmark(_GoBack) / ". Keith was lucky..." ->
#           merge into a single run " her window. Keith was lucky..." and
#           remove the _GoBack bookmark from this spot (it moves elsewhere
#           below, mirroring how Word relocates _GoBack to the most recent
#           edit point).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$d.Content.Find.Execute(". Keith was lucky", $true, $false, $false, $false, `
    $false, $true, 1, $false, ". Keith was lucky", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: split "...lets-slip-into-the-bathroom kind of guy..." so that
#           "-for-a-kinky-fuck" is inserted as its own run between
#           "...bathroom" and " kind of guy...", then stamp a fresh
#           _GoBack bookmark at the very end of that paragraph.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$splitPos = $full.IndexOf("lets-slip-into-the-bathroom kind") + "lets-slip-into-the-bathroom".Length

# Insert a unique placeholder right at the split point, then isolate it with
# a temporary bookmark so the surrounding same-formatted runs don't get
# reflowed into one another; replace the placeholder text from inside that
# bookmark's own range so the new text lands in its own run.
$insertPoint = $d.Range($splitPos, $splitPos)
$insertPoint.InsertAfter("@@SPLIT@@")

$full = $d.Content.Text
$markIdx = $full.IndexOf("@@SPLIT@@")
$markRange = $d.Range($markIdx, $markIdx + 9)
$tempMark = $d.Bookmarks.Add("zTmpSplitMark", $markRange)

$tempMark.Range.Find.Execute("@@SPLIT@@", $true, $false, $false, $false, `
    $false, $true, 1, $false, "-for-a-kinky-fuck", 2) | Out-Null

$d.Bookmarks("zTmpSplitMark").Delete()

# Now place the _GoBack bookmark at the end of this paragraph (right after
# "...without his glasses to boot?"). A zero-width range landing exactly on
# the paragraph mark does not behave reliably, so insert a throwaway
# placeholder character, bookmark around it, then delete just the character
# text, leaving the bookmark collapsed in the right spot.
$full = $d.Content.Text
$paraEnd = $full.IndexOf("without his glasses to boot?") + "without his glasses to boot?".Length

$endInsert = $d.Range($paraEnd, $paraEnd)
$endInsert.InsertAfter("@")

$endWrap = $d.Range($paraEnd, $paraEnd + 1)
$newGoBack = $d.Bookmarks.Add("_GoBack", $endWrap)
$newGoBack.Range.Text = ""
